# Update the cryptos price/volume table with the latest scraped values.
# Note: some "Price" (column D) values look like plain numbers (e.g. 593.52),
# but the sheet stores them as text (Excel's regional "." thousands-separator
# quirk elsewhere in the column, e.g. "61.972.37", means the whole column must
# stay text). A leading apostrophe forces Excel to keep the entry as text,
# exactly like typing '593.52 into a cell, instead of letting COM coerce it
# into a floating point Double.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.045.57'
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').Value = '3.004.98'
$ws.Range('E3').Value = '  -0.26%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''593.52'
$ws.Range('E5').Value = '  +1.45%  '
$ws.Range('D6').Value = '''146.91'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D8').Value = '3.002.70'
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('E9').Value = '  -2.08%  '
$ws.Range('E10').Value = '  +6.77%  '
$ws.Range('D11').Value = '''0.149'
$ws.Range('E11').Value = '  -0.31%  '
$ws.Range('E12').Value = '  -0.71%  '
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('D14').Value = '''34.31'
$ws.Range('E14').Value = '  -1.29%  '
$ws.Range('E15').Value = '  +2.76%  '
$ws.Range('D16').Value = '3.499.17'
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D17').Value = '''6.98'
$ws.Range('E17').Value = '  -1.64%  '
$ws.Range('D18').Value = '62.015.14'
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('D19').Value = '3.004.27'
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('D20').Value = '''446.58'
$ws.Range('E20').Value = '  -2.95%  '
$ws.Range('D21').Value = '''14.11'
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('E22').Value = '  -0.63%  '
$ws.Range('D23').Value = '''7.37'
$ws.Range('E23').Value = '  -0.95%  '
$ws.Range('D24').Value = '''82.16'
$ws.Range('E24').Value = '  +0.58%  '
$ws.Range('D25').Value = '''11.01'
$ws.Range('E25').Value = '  +9.24%  '
$ws.Range('D26').Value = '''2.23'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').Value = '''12.13'
$ws.Range('E27').Value = '  -2.08%  '
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('E29').Value = '  +2.98%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('E31').Value = '  +2.19%  '
$ws.Range('E32').Value = '  -0.51%  '
$ws.Range('E33').Value = '  -2.50%  '
$ws.Range('D34').Value = '''0.110'
$ws.Range('E34').Value = '  +1.04%  '
$ws.Range('D35').Value = '0.0₃0847'
$ws.Range('E35').Value = '  +4.07%  '
$ws.Range('E36').Value = '  -0.49%  '
$ws.Range('D37').Value = '''5.81'
$ws.Range('E37').Value = '  +0.48%  '
$ws.Range('D38').Value = '''50.16'
$ws.Range('E38').Value = '  -0.50%  '
$ws.Range('E39').Value = '  -3.96%  '
$ws.Range('D40').Value = '''8.98'
$ws.Range('E40').Value = '  -1.75%  '
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('E42').Value = '  +3.42%  '
$ws.Range('D43').Value = '''41.04'
$ws.Range('E43').Value = '  +9.84%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').Value = '''0.279'
$ws.Range('E44').Value = '  +3.63%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').Value = '''393.43'
$ws.Range('E45').Value = '  +0.77%  '
$ws.Range('E46').Value = '  -2.36%  '
$ws.Range('D47').Value = '2.715.02'
$ws.Range('E47').Value = '  -0.83%  '
$ws.Range('D48').Value = '''132.62'
$ws.Range('E48').Value = '  +2.38%  '
$ws.Range('D50').Value = '''2.17'
$ws.Range('E51').Value = '  -1.74%  '
